$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (single decimal point, pure digits) -- force them to stay Text like the
# original inline-string cells, then drop the format delta so no stray style
# sticks to the cell.
$textCells = @(
    'D5',
    'D6',
    'D7',
    'D9',
    'D10',
    'D12',
    'D14',
    'D15',
    'D17',
    'D18',
    'D20',
    'D21',
    'D23',
    'D24',
    'D25',
    'D26',
    'D27',
    'D29',
    'D32',
    'D33',
    'D35',
    'D36',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D45',
    'D46',
    'D47',
    'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '51.825.21'
$ws.Range('E2').Value = '  +0.43%  '
# Row 3
$ws.Range('E3').Value = '  +3.93%  '
# Row 4
$ws.Range('E4').Value = '  +0.00%  '
# Row 5
$ws.Range('D5').Value = '388.58'
$ws.Range('E5').Value = '  +1.71%  '
# Row 6
$ws.Range('D6').Value = '103.61'
$ws.Range('E6').Value = '  -0.80%  '
# Row 7
$ws.Range('D7').Value = '0.544'
$ws.Range('E7').Value = '  -0.34%  '
# Row 8
$ws.Range('E8').Value = '  -0.02%  '
# Row 9
$ws.Range('D9').Value = '0.590'
$ws.Range('E9').Value = '  -1.44%  '
# Row 10
$ws.Range('D10').Value = '37.21'
$ws.Range('E10').Value = '  +1.24%  '
# Row 12
$ws.Range('D12').Value = '0.0863'
$ws.Range('E12').Value = '  +0.30%  '
# Row 13
$ws.Range('D13').Value = '3.595.30'
$ws.Range('E13').Value = '  +3.97%  '
# Row 14
$ws.Range('D14').Value = '18.72'
$ws.Range('E14').Value = '  +0.95%  '
# Row 15
$ws.Range('D15').Value = '7.92'
$ws.Range('E15').Value = '  +0.53%  '
# Row 16
$ws.Range('D16').Value = '3.113.11'
$ws.Range('E16').Value = '  +3.60%  '
# Row 17
$ws.Range('D17').Value = '0.986'
$ws.Range('E17').Value = '  -0.98%  '
# Row 18
$ws.Range('D18').Value = '10.66'
$ws.Range('E18').Value = '  -4.61%  '
# Row 19
$ws.Range('D19').Value = '51.894.23'
$ws.Range('E19').Value = '  +0.52%  '
# Row 20
$ws.Range('D20').Value = '3.21'
$ws.Range('E20').Value = '  +3.33%  '
# Row 21
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').Value = '12.49'
$ws.Range('E21').Value = '  -0.93%  '
# Row 22
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').Value = '0.0₃0970'
$ws.Range('E22').Value = '  +0.67%  '
# Row 23
$ws.Range('D23').Value = '70.06'
$ws.Range('E23').Value = '  -0.68%  '
# Row 24
$ws.Range('D24').Value = '268.63'
$ws.Range('E24').Value = '  +0.29%  '
# Row 25
$ws.Range('D25').Value = '3.13'
$ws.Range('E25').Value = '  -3.18%  '
# Row 26
$ws.Range('D26').Value = '8.10'
$ws.Range('E26').Value = '  +1.96%  '
# Row 27
$ws.Range('D27').Value = '27.14'
$ws.Range('E27').Value = '  +3.82%  '
# Row 28
$ws.Range('E28').Value = '  -0.03%  '
# Row 29
$ws.Range('D29').Value = '7.19'
$ws.Range('E29').Value = '  -0.94%  '
# Row 30
$ws.Range('E30').Value = '  +0.06%  '
# Row 31
$ws.Range('E31').Value = '  -0.64%  '
# Row 32
$ws.Range('D32').Value = '10.37'
$ws.Range('E32').Value = '  -0.29%  '
# Row 33
$ws.Range('D33').Value = '35.48'
$ws.Range('E33').Value = '  +2.35%  '
# Row 34
$ws.Range('E34').Value = '  +1.73%  '
# Row 35
$ws.Range('D35').Value = '50.48'
$ws.Range('E35').Value = '  -1.86%  '
# Row 36
$ws.Range('D36').Value = '0.0450'
$ws.Range('E36').Value = '  +0.59%  '
# Row 37
$ws.Range('E37').Value = '  -0.18%  '
# Row 38
$ws.Range('E38').Value = '  +3.90%  '
# Row 39
$ws.Range('D39').Value = '0.290'
$ws.Range('E39').Value = '  +6.85%  '
# Row 40
$ws.Range('D40').Value = '1.89'
$ws.Range('E40').Value = '  +2.46%  '
# Row 41
$ws.Range('D41').Value = '16.90'
$ws.Range('E41').Value = '  -0.47%  '
# Row 42
$ws.Range('D42').Value = '2.60'
$ws.Range('E42').Value = '  +0.13%  '
# Row 43
$ws.Range('D43').Value = '128.23'
$ws.Range('E43').Value = '  +4.72%  '
# Row 44
$ws.Range('E44').Value = '  -0.33%  '
# Row 45
$ws.Range('D45').Value = '3.70'
$ws.Range('E45').Value = '  -4.33%  '
# Row 46
$ws.Range('D46').Value = '22.32'
$ws.Range('E46').Value = '  +4.58%  '
# Row 47
$ws.Range('D47').Value = '2.51'
$ws.Range('E47').Value = '  +6.65%  '
# Row 48
$ws.Range('E48').Value = '  +2.13%  '
# Row 49
$ws.Range('D49').Value = '2.048.23'
$ws.Range('E49').Value = '  +0.43%  '
# Row 50
$ws.Range('D50').Value = '3.415.93'
$ws.Range('E50').Value = '  +3.91%  '
# Row 51
$ws.Range('D51').Value = '0.207'
$ws.Range('E51').Value = '  +5.93%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
